# Auto-generated script applying value updates from the commit diff.
# Values correspond to market price / leve profit recalculations across
# the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 29382.273
$ws.Range("J87").Value = 29382.273
$ws.Range("L87").Value = 29382.273
$ws.Range("N87").Value = -31878.273
$ws.Range("H90").Value = 29382.273
$ws.Range("J90").Value = 29382.273
$ws.Range("L90").Value = 88146.819
$ws.Range("N90").Value = -100626.819
$ws.Range("H129").Value = 4630701
$ws.Range("J129").Value = 996.96155
$ws.Range("L129").Value = 2990.88465
$ws.Range("N129").Value = -12990.88465
$ws.Range("H137").Value = 1890019.9
$ws.Range("I137").Value = 4004584.8
$ws.Range("J137").Value = 2015.6786
$ws.Range("K137").Value = 12013754.4
$ws.Range("L137").Value = 6047.0358
$ws.Range("M137").Value = -12011204.4
$ws.Range("N137").Value = -11147.0358

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1234.6666
$ws.Range("I45").Value = 1012.6977
$ws.Range("J45").Value = 6007
$ws.Range("K45").Value = 1012.6977
$ws.Range("L45").Value = 6007
$ws.Range("M45").Value = -635.6977000000001
$ws.Range("N45").Value = -6761
$ws.Range("H135").Value = 25112.791
$ws.Range("J135").Value = 25112.791
$ws.Range("L135").Value = 25112.791
$ws.Range("N135").Value = -35252.791

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H101").Value = 3000
$ws.Range("J101").Value = 3000
$ws.Range("L101").Value = 9000
$ws.Range("N101").Value = -13868
$ws.Range("H105").Value = 2438.6667
$ws.Range("I105").Value = 2005
$ws.Range("J105").Value = 2655.5
$ws.Range("K105").Value = 2005
$ws.Range("L105").Value = 2655.5
$ws.Range("M105").Value = -258
$ws.Range("N105").Value = -6149.5
$ws.Range("H132").Value = 20220.64
$ws.Range("J132").Value = 20220.64
$ws.Range("L132").Value = 20220.64
$ws.Range("N132").Value = -30340.64
$ws.Range("H134").Value = 2882.1365
$ws.Range("I134").Value = 2755.342
$ws.Range("J134").Value = 3685.1667
$ws.Range("K134").Value = 8266.026
$ws.Range("L134").Value = 11055.5001
$ws.Range("M134").Value = -5731.026
$ws.Range("N134").Value = -16125.5001
$ws.Range("H135").Value = 23120.453
$ws.Range("I135").Value = 29000
$ws.Range("J135").Value = 23041
$ws.Range("K135").Value = 29000
$ws.Range("L135").Value = 23041
$ws.Range("M135").Value = -23930
$ws.Range("N135").Value = -33181
$ws.Range("H138").Value = 29000
$ws.Range("J138").Value = 29000
$ws.Range("L138").Value = 29000
$ws.Range("N138").Value = -39280
$ws.Range("H141").Value = 35884.617
$ws.Range("J141").Value = 29611.111
$ws.Range("L141").Value = 29611.111
$ws.Range("N141").Value = -39971.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2158.5
$ws.Range("I22").Value = 483.66666
$ws.Range("J22").Value = 3833.3333
$ws.Range("K22").Value = 483.66666
$ws.Range("L22").Value = 3833.3333
$ws.Range("M22").Value = -133.66666
$ws.Range("N22").Value = -4533.3333
$ws.Range("H132").Value = 3260.2632
$ws.Range("I132").Value = 2972.5386
$ws.Range("J132").Value = 3883.6667
$ws.Range("K132").Value = 8917.6158
$ws.Range("L132").Value = 11651.0001
$ws.Range("M132").Value = -6387.6158
$ws.Range("N132").Value = -16711.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 947.2059
$ws.Range("I5").Value = 613.0645
$ws.Range("J5").Value = 4400
$ws.Range("K5").Value = 1839.1935
$ws.Range("L5").Value = 13200
$ws.Range("M5").Value = -1727.1935
$ws.Range("N5").Value = -13424
$ws.Range("H102").Value = 2928.5715
$ws.Range("I102").Value = 2000
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 6000
$ws.Range("L102").Value = 9000
$ws.Range("M102").Value = -3566
$ws.Range("N102").Value = -13868
$ws.Range("H104").Value = 2984.1667
$ws.Range("J104").Value = 2984.1667
$ws.Range("L104").Value = 8952.500100000001
$ws.Range("N104").Value = -14194.5001
$ws.Range("H132").Value = 2657.6667
$ws.Range("I132").Value = 1449.8
$ws.Range("K132").Value = 13048.2
$ws.Range("M132").Value = -10518.2
$ws.Range("H135").Value = 947.2059
$ws.Range("I135").Value = 613.0645
$ws.Range("J135").Value = 4400
$ws.Range("K135").Value = 5517.5805
$ws.Range("L135").Value = 39600
$ws.Range("M135").Value = -2982.5805
$ws.Range("N135").Value = -44670
$ws.Range("H140").Value = 7652.5835
$ws.Range("I140").Value = 1590
$ws.Range("J140").Value = 9984.346
$ws.Range("K140").Value = 4770
$ws.Range("L140").Value = 29953.038
$ws.Range("M140").Value = 410
$ws.Range("N140").Value = -40313.038

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32080

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2116.6667
$ws.Range("J46").Value = 5400
$ws.Range("L46").Value = 5400
$ws.Range("N46").Value = -5776
$ws.Range("H68").Value = 5398.5713
$ws.Range("I68").Value = 5000
$ws.Range("J68").Value = 5465
$ws.Range("K68").Value = 5000
$ws.Range("L68").Value = 5465
$ws.Range("M68").Value = -4251
$ws.Range("N68").Value = -6963
$ws.Range("H71").Value = 5398.5713
$ws.Range("I71").Value = 5000
$ws.Range("J71").Value = 5465
$ws.Range("K71").Value = 25000
$ws.Range("L71").Value = 27325
$ws.Range("M71").Value = -21256
$ws.Range("N71").Value = -34813
$ws.Range("H132").Value = 3872.4443
$ws.Range("I132").Value = 3533.1667
$ws.Range("J132").Value = 4042.0833
$ws.Range("K132").Value = 10599.5001
$ws.Range("L132").Value = 12126.2499
$ws.Range("M132").Value = -8069.500100000001
$ws.Range("N132").Value = -17186.2499
$ws.Range("H134").Value = 16994.186
$ws.Range("J134").Value = 16994.186
$ws.Range("L134").Value = 16994.186
$ws.Range("N134").Value = -27134.186
$ws.Range("H136").Value = 6675334
$ws.Range("I136").Value = 11118778
$ws.Range("J136").Value = 10167.5
$ws.Range("K136").Value = 33356334
$ws.Range("L136").Value = 30502.5
$ws.Range("M136").Value = -33353784
$ws.Range("N136").Value = -35602.5
$ws.Range("H137").Value = 20001
$ws.Range("J137").Value = 20001
$ws.Range("L137").Value = 20001
$ws.Range("N137").Value = -30201

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1714
$ws.Range("I81").Value = 1070.909
$ws.Range("J81").Value = 2500
$ws.Range("K81").Value = 2141.818
$ws.Range("L81").Value = 5000
$ws.Range("M81").Value = -1080.818
$ws.Range("N81").Value = -7122
$ws.Range("H84").Value = 1714
$ws.Range("I84").Value = 1070.909
$ws.Range("J84").Value = 2500
$ws.Range("K84").Value = 10709.09
$ws.Range("L84").Value = 25000
$ws.Range("M84").Value = -5405.09
$ws.Range("N84").Value = -35608
$ws.Range("H113").Value = 1582.8572
$ws.Range("I113").Value = 156
$ws.Range("K113").Value = 468
$ws.Range("M113").Value = 1702
$ws.Range("H135").Value = 20467.959
$ws.Range("J135").Value = 20467.959
$ws.Range("L135").Value = 20467.959
$ws.Range("N135").Value = -30607.959
$ws.Range("H136").Value = 2769.1177
$ws.Range("I136").Value = 2119.3333
$ws.Range("K136").Value = 6357.999899999999
$ws.Range("M136").Value = -3807.999899999999
$ws.Range("H139").Value = 18955.8
$ws.Range("J139").Value = 18955.8
$ws.Range("L139").Value = 18955.8
$ws.Range("N139").Value = -29235.8
$ws.Range("H141").Value = 19756.512
$ws.Range("J141").Value = 19756.512
$ws.Range("L141").Value = 19756.512
$ws.Range("N141").Value = -30116.512
